$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (not ambiguous with numbers)
$ws.Range('D2').Value = '66.941.64'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '3.079.52'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('E6').Value = '  -2.39%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.076.07'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').Value = '3.588.63'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '66.893.67'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').Value = '3.078.61'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('E21').Value = '  +3.32%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('E25').Value = '  -4.12%  '
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('E27').Value = '  +3.70%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  -4.97%  '
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('E32').Value = '  -2.53%  '
$ws.Range('E33').Value = '  -2.21%  '
$ws.Range('D34').Value = '0.0₃0909'
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('E40').Value = '  -4.01%  '
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('D43').Value = '2.778.24'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E45').Value = '  -2.53%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('E51').Value = '  -1.25%  '

# Cells whose new values look numeric (e.g. "1.00", "7.70") must be forced to text
# so they are stored the same way as the original inline-string values.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '577.13'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '168.09'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.515'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.40'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.472'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000242'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.11'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.02'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.67'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '491.60'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.70'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.688'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.70'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.87'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.22'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.20'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.89'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.69'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.956'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '46.63'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.99'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.303'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.33'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '369.98'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0344'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '135.53'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.47'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '24.74'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.16'
